$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 102350
$ws.Cells.Item(3, 2).Value = 88955
$ws.Cells.Item(4, 2).Value = 45654
$ws.Cells.Item(5, 2).Value = 44008
$ws.Cells.Item(6, 2).Value = 42053
$ws.Cells.Item(7, 2).Value = 36194
$ws.Cells.Item(9, 2).Value = 27616
$ws.Cells.Item(10, 2).Value = 26724
$ws.Cells.Item(11, 2).Value = 26280
$ws.Cells.Item(12, 2).Value = 25729
$ws.Cells.Item(13, 2).Value = 25555
$ws.Cells.Item(14, 2).Value = 22486
$ws.Cells.Item(15, 2).Value = 21915
$ws.Cells.Item(16, 2).Value = 21857
$ws.Cells.Item(17, 2).Value = 21135
$ws.Cells.Item(18, 2).Value = 20853
$ws.Cells.Item(19, 2).Value = 18174
$ws.Cells.Item(20, 2).Value = 18083
$ws.Cells.Item(21, 2).Value = 17107
$ws.Cells.Item(22, 2).Value = 16471
$ws.Cells.Item(23, 2).Value = 15521
$ws.Cells.Item(24, 2).Value = 15220
$ws.Cells.Item(25, 2).Value = 15149
$ws.Cells.Item(26, 2).Value = 14653
$ws.Cells.Item(27, 2).Value = 14357
$ws.Cells.Item(28, 2).Value = 13996
$ws.Cells.Item(29, 2).Value = 13980
$ws.Cells.Item(31, 2).Value = 12803
$ws.Cells.Item(32, 2).Value = 12550
$ws.Cells.Item(33, 2).Value = 12494
$ws.Cells.Item(34, 2).Value = 12240
$ws.Cells.Item(35, 2).Value = 11968
$ws.Cells.Item(36, 2).Value = 11724
$ws.Cells.Item(37, 2).Value = 11479
$ws.Cells.Item(38, 2).Value = 10956
$ws.Cells.Item(39, 2).Value = 10710
$ws.Cells.Item(40, 2).Value = 10651
$ws.Cells.Item(41, 2).Value = 10537
$ws.Cells.Item(43, 2).Value = 10351
$ws.Cells.Item(44, 2).Value = 10283
$ws.Cells.Item(45, 2).Value = 10171
$ws.Cells.Item(46, 2).Value = 10090
$ws.Cells.Item(47, 2).Value = 10014
$ws.Cells.Item(48, 2).Value = 9903
$ws.Cells.Item(49, 2).Value = 9741
$ws.Cells.Item(50, 2).Value = 9726
$ws.Cells.Item(51, 2).Value = 9664
$ws.Cells.Item(52, 2).Value = 9640
$ws.Cells.Item(53, 2).Value = 9585
$ws.Cells.Item(54, 2).Value = 9510
$ws.Cells.Item(55, 2).Value = 9494
$ws.Cells.Item(56, 2).Value = 9413
$ws.Cells.Item(57, 2).Value = 9105
$ws.Cells.Item(58, 2).Value = 9013
$ws.Cells.Item(59, 2).Value = 8966
$ws.Cells.Item(60, 2).Value = 8754
$ws.Cells.Item(61, 2).Value = 8493
$ws.Cells.Item(62, 2).Value = 8211
$ws.Cells.Item(63, 2).Value = 8209
$ws.Cells.Item(64, 2).Value = 8099
$ws.Cells.Item(65, 2).Value = 7996
$ws.Cells.Item(66, 2).Value = 7910
$ws.Cells.Item(67, 2).Value = 7867
$ws.Cells.Item(68, 2).Value = 7813
$ws.Cells.Item(69, 2).Value = 7818
$ws.Cells.Item(70, 2).Value = 7765
$ws.Cells.Item(71, 2).Value = 7746
$ws.Cells.Item(72, 2).Value = 7703
$ws.Cells.Item(73, 2).Value = 7632
$ws.Cells.Item(74, 2).Value = 7523
$ws.Cells.Item(75, 2).Value = 7456
$ws.Cells.Item(76, 2).Value = 7434
$ws.Cells.Item(77, 2).Value = 7439
$ws.Cells.Item(78, 2).Value = 7368
$ws.Cells.Item(79, 2).Value = 7285
$ws.Cells.Item(81, 2).Value = 7156
$ws.Cells.Item(82, 2).Value = 7112
$ws.Cells.Item(83, 2).Value = 7075
$ws.Cells.Item(84, 2).Value = 7061
$ws.Cells.Item(85, 2).Value = 7049
$ws.Cells.Item(86, 2).Value = 7027
$ws.Cells.Item(87, 2).Value = 7027
$ws.Cells.Item(88, 2).Value = 7015
$ws.Cells.Item(89, 2).Value = 6971
$ws.Cells.Item(90, 2).Value = 6963
$ws.Cells.Item(91, 2).Value = 6956
$ws.Cells.Item(92, 2).Value = 6929
$ws.Cells.Item(93, 2).Value = 6851
$ws.Cells.Item(94, 2).Value = 6833
$ws.Cells.Item(95, 2).Value = 6817
$ws.Cells.Item(96, 2).Value = 6794
$ws.Cells.Item(97, 2).Value = 6725
$ws.Cells.Item(98, 2).Value = 6705
